$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": M6 15693.84 -> 25219.13
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M6").Value = 25219.13

# Sheet "VENTA MENSUAL": F6 15693.84 -> 25219.13 ; F19 22179.48 -> 31704.77
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F6").Value = 25219.13
$ws2.Range("F19").Value = 31704.77

# Sheet "CUMPLIMIENTO MENSUAL": row 16 and row 19 updates
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 29659.77
$ws3.Range("E16").Value = -1449.93
$ws3.Range("F16").Value = 1.05139802281757

$ws3.Range("D19").Value = 31704.77
$ws3.Range("E19").Value = 15514.53386304603
$ws3.Range("F19").Value = 0.6714366245626133
